# "La danse du lion" -> "La dance du lion", with the run split into
# "La " / "dance" / " du lion" and the _GoBack bookmark re-anchored right
# after "dance" (mirrors how Word marks the last-edit location after a
# manual retype of "danse" -> "dance").
$d = $word.ActiveDocument

# Replace "danse" (chars 4-8, 1-based range 3..8) with "dance".
$r = $d.Range(3, 8)
$r.Text = "dance"

# Force a run boundary right after "La " (position 3) by dropping a
# transient bookmark there, then move the document's _GoBack bookmark to
# sit immediately after "dance" (position 8) -- re-adding a bookmark with
# the same name relocates it instead of creating a duplicate.
$null = $d.Bookmarks.Add("_GoBack", $d.Range(3, 3))
$null = $d.Bookmarks.Add("_GoBack", $d.Range(8, 8))
